$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70's date cell (A70) is no longer the last row, so it switches from the
# date-only format to the regular date+time format used by the rest of the rows.
$ws.Range("A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new row 71 with the day's values.
$ws.Range("A71").Value = 45811
$ws.Range("B71").Value = 303
$ws.Range("C71").Value = 300
$ws.Range("D71").Value = 303

# The new last row (A71) gets the date-only format.
$ws.Range("A71").NumberFormat = "YYYY-MM-DD"
